$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Gnai2"
$ws.Cells.Item(2, 3).Value = "Agtr2"
$ws.Cells.Item(2, 4).Value = "FAPs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 150.0354306666667
$ws.Cells.Item(2, 8).Value = 450.106292
$ws.Cells.Item(2, 9).Value = 0.4152507364956075
$ws.Cells.Item(2, 10).Value = 0.4152507364956075
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 34.682839
$ws.Cells.Item(2, 14).Value = 104.048517
$ws.Cells.Item(2, 15).Value = 0.9919940127517238
$ws.Cells.Item(2, 16).Value = 0.9919940127517238
$ws.Cells.Item(2, 17).Value = 5203.654686107662
$ws.Cells.Item(2, 18).Value = 46832.89217496897
$ws.Cells.Item(2, 19).Value = 0.4119262443943864
$ws.Cells.Item(2, 20).Value = 0.4119262443943864

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Gnai2"
$ws.Cells.Item(3, 3).Value = "Agtr2"
$ws.Cells.Item(3, 4).Value = "sCs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 150.0354306666667
$ws.Cells.Item(3, 8).Value = 450.106292
$ws.Cells.Item(3, 9).Value = 0.4152507364956075
$ws.Cells.Item(3, 10).Value = 0.4152507364956075
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.2799113333333333
$ws.Cells.Item(3, 14).Value = 0.839734
$ws.Cells.Item(3, 15).Value = 0.008005987248276263
$ws.Cells.Item(3, 16).Value = 0.008005987248276263
$ws.Cells.Item(3, 17).Value = 41.99661744514756
$ws.Cells.Item(3, 18).Value = 377.969557006328
$ws.Cells.Item(3, 19).Value = 0.00332449210122116
$ws.Cells.Item(3, 20).Value = 0.00332449210122116

# Row 4
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Gnai2"
$ws.Cells.Item(4, 3).Value = "Agtr2"
$ws.Cells.Item(4, 4).Value = "FAPs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 68.382243
$ws.Cells.Item(4, 8).Value = 205.146729
$ws.Cells.Item(4, 9).Value = 0.1892604742946246
$ws.Cells.Item(4, 10).Value = 0.1892604742946246
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 34.682839
$ws.Cells.Item(4, 14).Value = 104.048517
$ws.Cells.Item(4, 15).Value = 0.9919940127517238
$ws.Cells.Item(4, 16).Value = 0.9919940127517238
$ws.Cells.Item(4, 17).Value = 2371.690324427877
$ws.Cells.Item(4, 18).Value = 21345.21291985089
$ws.Cells.Item(4, 19).Value = 0.1877452573508192
$ws.Cells.Item(4, 20).Value = 0.1877452573508191

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Gnai2"
$ws.Cells.Item(5, 3).Value = "Agtr2"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 68.382243
$ws.Cells.Item(5, 8).Value = 205.146729
$ws.Cells.Item(5, 9).Value = 0.1892604742946246
$ws.Cells.Item(5, 10).Value = 0.1892604742946246
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.2799113333333333
$ws.Cells.Item(5, 14).Value = 0.839734
$ws.Cells.Item(5, 15).Value = 0.008005987248276263
$ws.Cells.Item(5, 16).Value = 0.008005987248276263
$ws.Cells.Item(5, 17).Value = 19.140964814454
$ws.Cells.Item(5, 18).Value = 172.268683330086
$ws.Cells.Item(5, 19).Value = 0.001515216943805482
$ws.Cells.Item(5, 20).Value = 0.001515216943805482

# Row 6
$ws.Cells.Item(6, 1).Value = "M2"
$ws.Cells.Item(6, 2).Value = "Gnai2"
$ws.Cells.Item(6, 3).Value = "Agtr2"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 104.737245
$ws.Cells.Item(6, 8).Value = 314.211735
$ws.Cells.Item(6, 9).Value = 0.2898796499701289
$ws.Cells.Item(6, 10).Value = 0.2898796499701289
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 34.682839
$ws.Cells.Item(6, 14).Value = 104.048517
$ws.Cells.Item(6, 15).Value = 0.9919940127517238
$ws.Cells.Item(6, 16).Value = 0.9919940127517238
$ws.Cells.Item(6, 17).Value = 3632.585005638555
$ws.Cells.Item(6, 18).Value = 32693.265050747
$ws.Cells.Item(6, 19).Value = 0.2875588771889333
$ws.Cells.Item(6, 20).Value = 0.2875588771889332

# Row 7
$ws.Cells.Item(7, 1).Value = "M2"
$ws.Cells.Item(7, 2).Value = "Gnai2"
$ws.Cells.Item(7, 3).Value = "Agtr2"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 104.737245
$ws.Cells.Item(7, 8).Value = 314.211735
$ws.Cells.Item(7, 9).Value = 0.2898796499701289
$ws.Cells.Item(7, 10).Value = 0.2898796499701289
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.2799113333333333
$ws.Cells.Item(7, 14).Value = 0.839734
$ws.Cells.Item(7, 15).Value = 0.008005987248276263
$ws.Cells.Item(7, 16).Value = 0.008005987248276263
$ws.Cells.Item(7, 17).Value = 29.31714189761
$ws.Cells.Item(7, 18).Value = 263.85427707849
$ws.Cells.Item(7, 19).Value = 0.002320772781195639
$ws.Cells.Item(7, 20).Value = 0.002320772781195639

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Gnai2"
$ws.Cells.Item(8, 3).Value = "Agtr2"
$ws.Cells.Item(8, 4).Value = "FAPs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 38.15794
$ws.Cells.Item(8, 8).Value = 114.47382
$ws.Cells.Item(8, 9).Value = 0.105609139239639
$ws.Cells.Item(8, 10).Value = 0.105609139239639
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 34.682839
$ws.Cells.Item(8, 14).Value = 104.048517
$ws.Cells.Item(8, 15).Value = 0.9919940127517238
$ws.Cells.Item(8, 16).Value = 0.9919940127517238
$ws.Cells.Item(8, 17).Value = 1323.42568959166
$ws.Cells.Item(8, 18).Value = 11910.83120632494
$ws.Cells.Item(8, 19).Value = 0.104763633817585
$ws.Cells.Item(8, 20).Value = 0.104763633817585

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Gnai2"
$ws.Cells.Item(9, 3).Value = "Agtr2"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 38.15794
$ws.Cells.Item(9, 8).Value = 114.47382
$ws.Cells.Item(9, 9).Value = 0.105609139239639
$ws.Cells.Item(9, 10).Value = 0.105609139239639
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.2799113333333333
$ws.Cells.Item(9, 14).Value = 0.839734
$ws.Cells.Item(9, 15).Value = 0.008005987248276263
$ws.Cells.Item(9, 16).Value = 0.008005987248276263
$ws.Cells.Item(9, 17).Value = 10.68083986265333
$ws.Cells.Item(9, 18).Value = 96.12755876388
$ws.Cells.Item(9, 19).Value = 0.0008455054220539822
$ws.Cells.Item(9, 20).Value = 0.000845505422053982

